# Scheduled-runner refresh of market-board derived columns
# (currentAveragePrice[, NQ, HQ] / LevePrice[NQ, HQ] / LeveProfit[NQ, HQ])
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR Leve-profit sheets.
# Values below are the refreshed figures pulled in by the scheduled run;
# cells that no longer have a computed profit are cleared, and cells that
# newly gained one are populated.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 346
$ws.Range("J2").Value = 352.4
$ws.Range("L2").Value = 352.4
$ws.Range("N2").Value = -578.4

$ws.Range("H33").Value = 370.70587
$ws.Range("I33").Value = 338
$ws.Range("J33").Value = 430.66666
$ws.Range("K33").Value = 338
$ws.Range("L33").Value = 430.66666
$ws.Range("M33").Value = -109
$ws.Range("N33").Value = -888.66666

$ws.Range("H64").Value = 16317.333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 16317.333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 16317.333
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -16813.333

$ws.Range("H67").Value = 16317.333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 16317.333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 16317.333
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -18033.333

$ws.Range("H75").Value = 28000
$ws.Range("J75").Value = 28000
$ws.Range("L75").Value = 28000
$ws.Range("N75").Value = -29872

$ws.Range("H78").Value = 28000
$ws.Range("J78").Value = 28000
$ws.Range("L78").Value = 84000
$ws.Range("N78").Value = -93360

$ws.Range("H93").Value = 1000000
$ws.Range("J93").Value = 1000000
$ws.Range("L93").Value = 1000000
$ws.Range("N93").Value = -1004992

$ws.Range("H95").Value = 72268.836
$ws.Range("J95").Value = 72268.836
$ws.Range("L95").Value = 72268.836
$ws.Range("N95").Value = -77760.836

$ws.Range("H137").Value = 55560224
$ws.Range("I137").Value = 100003736
$ws.Range("J137").Value = 5837.875
$ws.Range("K137").Value = 300011208
$ws.Range("L137").Value = 17513.625
$ws.Range("M137").Value = -300008658
$ws.Range("N137").Value = -22613.625

$ws.Range("H138").Value = 5325.213
$ws.Range("J138").Value = 5266.7964
$ws.Range("L138").Value = 15800.3892
$ws.Range("N138").Value = -26080.3892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11093.472
$ws.Range("I32").Value = 8311.270500000001
$ws.Range("J32").Value = 37802.6
$ws.Range("K32").Value = 8311.270500000001
$ws.Range("L32").Value = 37802.6
$ws.Range("M32").Value = -8024.270500000001
$ws.Range("N32").Value = -38376.6

$ws.Range("H45").Value = 7504.5
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H61").Value = 3728.5264
$ws.Range("I61").Value = 2990.7144
$ws.Range("J61").Value = 12336.333
$ws.Range("K61").Value = 2990.7144
$ws.Range("L61").Value = 12336.333
$ws.Range("M61").Value = -2778.7144
$ws.Range("N61").Value = -12760.333

$ws.Range("H122").Value = 2663.825
$ws.Range("I122").Value = 2052.1667
$ws.Range("K122").Value = 6156.500100000001
$ws.Range("M122").Value = -3706.500100000001

$ws.Range("H136").Value = 3728.5264
$ws.Range("I136").Value = 2990.7144
$ws.Range("J136").Value = 12336.333
$ws.Range("K136").Value = 8972.143199999999
$ws.Range("L136").Value = 37008.999
$ws.Range("M136").Value = -6422.143199999999
$ws.Range("N136").Value = -42108.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3998.5
$ws.Range("I20").Value = 995
$ws.Range("J20").Value = 13009
$ws.Range("K20").Value = 995
$ws.Range("L20").Value = 13009
$ws.Range("M20").Value = -748
$ws.Range("N20").Value = -13503

$ws.Range("H107").Value = 1251.5
$ws.Range("I107").Value = 1251.5
$ws.Range("K107").Value = 1251.5
$ws.Range("M107").Value = 668.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 60028.8
$ws.Range("I31").Value = 11953.154
$ws.Range("K31").Value = 11953.154
$ws.Range("M31").Value = -11658.154

$ws.Range("H34").Value = 60028.8
$ws.Range("I34").Value = 11953.154
$ws.Range("K34").Value = 11953.154
$ws.Range("M34").Value = -11751.154

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H99").Value = 4417.8184
$ws.Range("I99").Value = 4616
$ws.Range("K99").Value = 4616
$ws.Range("M99").Value = -3118

$ws.Range("H126").Value = 4417.8184
$ws.Range("I126").Value = 4616
$ws.Range("K126").Value = 13848
$ws.Range("M126").Value = -11378

$ws.Range("H141").Value = 348066.75
$ws.Range("J141").Value = 364070.44
$ws.Range("L141").Value = 364070.44
$ws.Range("N141").Value = -374430.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1389981.8
$ws.Range("I5").Value = 997.94116
$ws.Range("J5").Value = 4763228
$ws.Range("K5").Value = 2993.82348
$ws.Range("L5").Value = 14289684
$ws.Range("M5").Value = -2881.82348
$ws.Range("N5").Value = -14289908

$ws.Range("H9").Value = 63493.125
$ws.Range("I9").Value = 316990
$ws.Range("J9").Value = 4993.846
$ws.Range("K9").Value = 950970
$ws.Range("L9").Value = 14981.538
$ws.Range("M9").Value = -950746
$ws.Range("N9").Value = -15429.538

$ws.Range("H38").Value = 57.615383
$ws.Range("I38").Value = 45.666668
$ws.Range("J38").Value = 67.85714
$ws.Range("K38").Value = 137.000004
$ws.Range("L38").Value = 203.57142
$ws.Range("M38").Value = 209.999996
$ws.Range("N38").Value = -897.57142

$ws.Range("H131").Value = 29042506
$ws.Range("J131").Value = 24307774
$ws.Range("L131").Value = 72923322
$ws.Range("N131").Value = -72933402

$ws.Range("H132").Value = 5435.9375
$ws.Range("J132").Value = 5997.3335
$ws.Range("L132").Value = 53976.0015
$ws.Range("N132").Value = -59036.0015

$ws.Range("H135").Value = 1389981.8
$ws.Range("I135").Value = 997.94116
$ws.Range("J135").Value = 4763228
$ws.Range("K135").Value = 8981.470439999999
$ws.Range("L135").Value = 42869052
$ws.Range("M135").Value = -6446.470439999999
$ws.Range("N135").Value = -42874122

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30302

$ws.Range("H46").Value = 35000
$ws.Range("I46").Value = 35000
$ws.Range("J46").Value = 35000
$ws.Range("K46").Value = 35000
$ws.Range("L46").Value = 35000
$ws.Range("M46").Value = -34844
$ws.Range("N46").Value = -35312

$ws.Range("H53").Value = 4444444
$ws.Range("J53").Value = 4444444
$ws.Range("L53").Value = 4444444
$ws.Range("N53").Value = -4445706

$ws.Range("H126").Value = 4026.6875
$ws.Range("I126").Value = 3983.9614
$ws.Range("K126").Value = 11951.8842
$ws.Range("M126").Value = -9481.8842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 142867780
$ws.Range("J22").Value = 14599.6
$ws.Range("L22").Value = 14599.6
$ws.Range("N22").Value = -15189.6

$ws.Range("H27").Value = 142867780
$ws.Range("J27").Value = 14599.6
$ws.Range("L27").Value = 14599.6
$ws.Range("N27").Value = -14813.6

$ws.Range("H123").Value = 72000
$ws.Range("J123").Value = 72000
$ws.Range("L123").Value = 72000
$ws.Range("N123").Value = -81800

$ws.Range("H130").Value = 74433
$ws.Range("J130").Value = 74433
$ws.Range("L130").Value = 74433
$ws.Range("N130").Value = -84473

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11315.3
$ws.Range("I81").Value = 10630.2
$ws.Range("J81").Value = 12000.4
$ws.Range("K81").Value = 21260.4
$ws.Range("L81").Value = 24000.8
$ws.Range("M81").Value = -20199.4
$ws.Range("N81").Value = -26122.8

$ws.Range("H84").Value = 11315.3
$ws.Range("I84").Value = 10630.2
$ws.Range("J84").Value = 12000.4
$ws.Range("K84").Value = 106302
$ws.Range("L84").Value = 120004
$ws.Range("M84").Value = -100998
$ws.Range("N84").Value = -130612

$ws.Range("H126").Value = 3311.1
$ws.Range("I126").Value = 1952.4166
$ws.Range("K126").Value = 5857.2498
$ws.Range("M126").Value = -3387.2498

$ws.Range("H132").Value = 3371.2666
$ws.Range("J132").Value = 5537.6665
$ws.Range("L132").Value = 16612.9995
$ws.Range("N132").Value = -21672.9995
